# Refresh Universalis market-price snapshots and recompute leve profit columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) per sheet.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 657.3333
$ws.Range("I28").Value = 398.75
$ws.Range("J28").Value = 1174.5
$ws.Range("K28").Value = 398.75
$ws.Range("L28").Value = 1174.5
$ws.Range("M28").Value = 86.25
$ws.Range("N28").Value = -2144.5
# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 3622.2
$ws.Range("I86").Value = 2223.1177
$ws.Range("J86").Value = 5451.769
$ws.Range("K86").Value = 2223.1177
$ws.Range("L86").Value = 5451.769
$ws.Range("M86").Value = -1100.1177
$ws.Range("N86").Value = -7697.769
# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 3622.2
$ws.Range("I89").Value = 2223.1177
$ws.Range("J89").Value = 5451.769
$ws.Range("K89").Value = 11115.5885
$ws.Range("L89").Value = 27258.845
$ws.Range("M89").Value = -5499.588499999998
$ws.Range("N89").Value = -38490.845

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 103: Sweeping the Legs
$ws.Range("H103").Value = 34350
$ws.Range("J103").Value = 34350
$ws.Range("L103").Value = 34350
$ws.Range("N103").Value = -36694

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 13751.467
$ws.Range("I82").Value = 2913.2
$ws.Range("J82").Value = 35428
$ws.Range("K82").Value = 2913.2
$ws.Range("L82").Value = 35428
$ws.Range("M82").Value = -2530.2
$ws.Range("N82").Value = -36194
# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 13751.467
$ws.Range("I85").Value = 2913.2
$ws.Range("J85").Value = 35428
$ws.Range("K85").Value = 2913.2
$ws.Range("L85").Value = 35428
$ws.Range("M85").Value = -1587.2
$ws.Range("N85").Value = -38080
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1775.7916
$ws.Range("I134").Value = 1839.9512
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 5519.8536
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -2984.8536
$ws.Range("N134").Value = -9270

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 21421.176
$ws.Range("J31").Value = 2361.5715
$ws.Range("L31").Value = 2361.5715
$ws.Range("N31").Value = -2951.5715
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 21421.176
$ws.Range("J34").Value = 2361.5715
$ws.Range("L34").Value = 2361.5715
$ws.Range("N34").Value = -2765.5715
# Row 103: Spare a Rod and Spoil the Fishers
$ws.Range("H103").Value = 12481.714
$ws.Range("I103").Value = 7474.4
$ws.Range("J103").Value = 25000
$ws.Range("K103").Value = 7474.4
$ws.Range("L103").Value = 25000
$ws.Range("M103").Value = -6302.4
$ws.Range("N103").Value = -27344
# Row 107: Built to Last
$ws.Range("H107").Value = 4188.4644
$ws.Range("I107").Value = 6545.294
$ws.Range("J107").Value = 546.0909
$ws.Range("K107").Value = 6545.294
$ws.Range("L107").Value = 546.0909
$ws.Range("M107").Value = -4625.294
$ws.Range("N107").Value = -4386.0909

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 125963.25
$ws.Range("I4").Value = 230855.23
$ws.Range("K4").Value = 692565.6900000001
$ws.Range("M4").Value = -692453.6900000001
# Row 5: What a Sap
$ws.Range("H5").Value = 7014.375
$ws.Range("I5").Value = 593.6957
$ws.Range("J5").Value = 23422.777
$ws.Range("K5").Value = 1781.0871
$ws.Range("L5").Value = 70268.33099999999
$ws.Range("M5").Value = -1669.0871
$ws.Range("N5").Value = -70492.33099999999
# Row 7: It's Always Sunny in Vylbrand
$ws.Range("H7").Value = 173
$ws.Range("I7").Value = 177.8
$ws.Range("J7").Value = 149
$ws.Range("K7").Value = 533.4000000000001
$ws.Range("L7").Value = 447
$ws.Range("M7").Value = -421.4000000000001
$ws.Range("N7").Value = -671
# Row 37: I Love Lamprey
$ws.Range("H37").Value = 458614.75
$ws.Range("J37").Value = 458614.75
$ws.Range("L37").Value = 1375844.25
$ws.Range("N37").Value = -1376068.25
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 1734.0441
$ws.Range("I68").Value = 843.6539
$ws.Range("K68").Value = 2530.9617
$ws.Range("M68").Value = -1719.9617
# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 1734.0441
$ws.Range("I71").Value = 843.6539
$ws.Range("K71").Value = 7592.8851
$ws.Range("M71").Value = -3536.8851
# Row 100: Souper
$ws.Range("H100").Value = 7499
$ws.Range("I100").Value = 545
$ws.Range("J100").Value = 9237.5
$ws.Range("K100").Value = 1635
$ws.Range("L100").Value = 27712.5
$ws.Range("M100").Value = -824
$ws.Range("N100").Value = -29334.5
# Row 107: Slippery Service
$ws.Range("H107").Value = 214591.92
$ws.Range("J107").Value = 480532.8
$ws.Range("L107").Value = 1441598.4
$ws.Range("N107").Value = -1445438.4
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1215.91
$ws.Range("I131").Value = 680
$ws.Range("J131").Value = 1238.2396
$ws.Range("K131").Value = 2040
$ws.Range("L131").Value = 3714.718800000001
$ws.Range("M131").Value = 3000
$ws.Range("N131").Value = -13794.7188
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 7014.375
$ws.Range("I135").Value = 593.6957
$ws.Range("J135").Value = 23422.777
$ws.Range("K135").Value = 5343.2613
$ws.Range("L135").Value = 210804.993
$ws.Range("M135").Value = -2808.2613
$ws.Range("N135").Value = -215874.993

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 531854.2
$ws.Range("J107").Value = 1010363.6
$ws.Range("L107").Value = 1010363.6
$ws.Range("N107").Value = -1014203.6
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2943.5454
$ws.Range("I113").Value = 3297.375
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 3297.375
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -1127.375
$ws.Range("N113").Value = -6340
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2799.524
$ws.Range("I122").Value = 2578.1333
$ws.Range("K122").Value = 7734.3999
$ws.Range("M122").Value = -5284.3999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 119: Fit for a Friend
$ws.Range("H119").Value = 39516
$ws.Range("J119").Value = 39516
$ws.Range("L119").Value = 39516
$ws.Range("N119").Value = -49192

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table
$ws.Range("H113").Value = 506
$ws.Range("I113").Value = 132
$ws.Range("J113").Value = 552.75
$ws.Range("K113").Value = 396
$ws.Range("L113").Value = 1658.25
$ws.Range("M113").Value = 1774
$ws.Range("N113").Value = -5998.25
# Row 119: A Job Well Done
$ws.Range("H119").Value = 16398.666
$ws.Range("J119").Value = 16398.666
$ws.Range("L119").Value = 16398.666
$ws.Range("N119").Value = -26074.666
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1526.421
$ws.Range("I126").Value = 1652.0769
$ws.Range("K126").Value = 4956.2307
$ws.Range("M126").Value = -2486.2307
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2965.5952
$ws.Range("I132").Value = 3255.8462
$ws.Range("K132").Value = 9767.5386
$ws.Range("M132").Value = -7237.5386
